$wb = $excel.ActiveWorkbook

$wsSolar = $wb.Worksheets.Item("SolarPV")
$wsSolar.Range("B2").Value = 5000
$wsSolar.Range("B2").Select()

$wsBattery = $wb.Worksheets.Item("Battery_MV")
$wsBattery.Range("C3").Select()
